$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the
# 9aa49f36-b2eb-4c8d-bc8c-18e31aac6966.md row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-23 22:43:00"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the 9aa49f36-b2eb-4c8d-bc8c-18e31aac6966.md row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-23 22:42:55"
$wsZhCn.Range("K4").Value = "2016-08-23 22:43:26"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the 9aa49f36-b2eb-4c8d-bc8c-18e31aac6966.md row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-08-23 22:43:00"
$wsDeDe.Range("K4").Value = "2016-08-23 22:43:34"
